$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net Rr"
$ws.Range("G1").Value = "Pts"

$ws.Range("A1:G10").Select()
